# Atualizando base de dados da pesquisa via Streamlit
# Adds the new survey response as row 17 at the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# Columns that carry real values in the new record.
$ws.Cells.Item($row, 1).Value = "teste@mrv.com.br"
$ws.Cells.Item($row, 9).Value = "2025-05-21 11:53:27"
$ws.Cells.Item($row, 21).Value = "PAP - Dossiê: Comentário teste 1; Painel Análises Forecast de Produção - PLNESROBR009: Comentário teste 2"
$ws.Cells.Item($row, 22).Value = "Planilha automatizada - teste,Objetivo 1,Excel,OUTROS,🪙 Importante,7.0; Ferramentas - Planejamento Operacional,Objetivo 2,Python,MOP/EMP,🟢 Pouco Importante,4.0"

# Remaining columns stay blank in the new record, same as every other row,
# but still need a real (empty-text) cell rather than a totally absent one.
# A plain "" assignment clears the cell outright, so seed it with a lone
# quote-prefix character (true-empty text in Excel) and then strip the
# quote-prefix formatting it implies, leaving a normal-styled empty cell.
$blankCols = 2,3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20
foreach ($c in $blankCols) {
    $cell = $ws.Cells.Item($row, $c)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
